$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.313208699226379
$ws.Range("B1").Value = 1.35558021068573
$ws.Range("C1").Value = 3.788367748260498
$ws.Range("D1").Value = 3.514906167984009
$ws.Range("E1").Value = 1.04081654548645
